# Apply the "Add more test cases and adjust data elements descriptions" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Adjust existing "Data Element" descriptions (column A) ---
$ws.Range("A4").Value = "Reference ID"
$ws.Range("A8").Value = "Post ID"

# --- New TC-4 column block (K:L), mirroring the TC-3 block in G:J ---

# Header "TC-4" merged across K1:L1, matching the style of the other
# TC-n headers (copy formatting from the TC-2 header so the existing
# blue/centered style gets reused instead of minting a new one).
$ws.Range("K1").Value = "TC-4"
$ws.Range("E1:F1").Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1:L1").Merge() | Out-Null

# Input / Expected Result sub-headers
$ws.Range("K2").Value = "Input"
$ws.Range("L2").Value = "Expected Result"

# Row 4: reference TC id
$ws.Range("K4").Value = "TC-1"

# Row 5: request name for the new test case
$ws.Range("K5").Value = "Put Comment"

# Row 8: numeric Input/Expected-Result pair (postID)
$ws.Range("K8").Value = 102
$ws.Range("L8").Formula = "=IF(K8<>"""",K8,"""")"

# Row 9-10: shared "title" formula block
$ws.Range("L9:L10").Formula = "=IF(K9<>"""",K9,"""")"
$ws.Range("K10").Value = "Test four"

# Row 11: "body" formula (mirrors the J-column body formula pattern)
$ws.Range("L11").Formula = "=IF(J11<>"""",J11,"""")"

# --- Selection moves to A7 ---
$ws.Range("A7").Select() | Out-Null
